$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K18").Value = 11500
$ws.Range("I18").Value = 11500
$ws.Range("M18").Value = -11216
$ws.Range("H18").Value = 11500
$ws.Range("K19").Value = 1221
$ws.Range("M19").Value = -1046
$ws.Range("I19").Value = 1221
$ws.Range("H19").Value = 1221
$ws.Range("I127").Value = 955.2
$ws.Range("K127").Value = 2865.6
$ws.Range("M127").Value = 2094.4
$ws.Range("N127").ClearContents()
$ws.Range("H127").Value = 955.2
$ws.Range("L127").Value = 0
$ws.Range("J127").Value = 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3363.3928
$ws.Range("M32").Value = -2830.6296
$ws.Range("I32").Value = 3117.6296
$ws.Range("K32").Value = 3117.6296
$ws.Range("H74").Value = 1252.8889
$ws.Range("K74").Value = 982.5714
$ws.Range("I74").Value = 982.5714
$ws.Range("M74").Value = -108.5714
$ws.Range("H77").Value = 1252.8889
$ws.Range("K77").Value = 4912.857
$ws.Range("M77").Value = -544.857
$ws.Range("I77").Value = 982.5714
$ws.Range("I97").Value = 599.5
$ws.Range("H97").Value = 599.5
$ws.Range("K97").Value = 599.5
$ws.Range("M97").Value = -103.5
$ws.Range("H122").Value = 1536.1538
$ws.Range("J122").Value = 1199.5
$ws.Range("L122").Value = 3598.5
$ws.Range("N122").Value = -8498.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2142.8235
$ws.Range("L20").Value = 4154
$ws.Range("N20").Value = -4648
$ws.Range("J20").Value = 4154
$ws.Range("I20").Value = 1304.8334
$ws.Range("M20").Value = -1057.8334
$ws.Range("K20").Value = 1304.8334
$ws.Range("H94").Value = 480.92856
$ws.Range("M94").Value = 31
$ws.Range("K94").Value = 420
$ws.Range("I94").Value = 420
$ws.Range("K107").Value = 553
$ws.Range("I107").Value = 553
$ws.Range("M107").Value = 1367
$ws.Range("H107").Value = 727.1429000000001
$ws.Range("M134").Value = -11335.6671
$ws.Range("I134").Value = 4623.5557
$ws.Range("K134").Value = 13870.6671
$ws.Range("H134").Value = 4361.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N12").ClearContents()
$ws.Range("L12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("L22").Value = 13333867
$ws.Range("K22").Value = 561.5714
$ws.Range("M22").Value = -211.5714
$ws.Range("I22").Value = 561.5714
$ws.Range("H22").Value = 4000553
$ws.Range("J22").Value = 13333867
$ws.Range("N22").Value = -13334567
$ws.Range("J31").Value = 1499.5
$ws.Range("H31").Value = 1430.091
$ws.Range("N31").Value = -2089.5
$ws.Range("L31").Value = 1499.5
$ws.Range("H34").Value = 1430.091
$ws.Range("J34").Value = 1499.5
$ws.Range("L34").Value = 1499.5
$ws.Range("N34").Value = -1903.5
$ws.Range("H62").Value = 3596.3333
$ws.Range("L62").Value = 3596.3333
$ws.Range("N62").Value = -4844.3333
$ws.Range("J62").Value = 3596.3333
$ws.Range("N65").Value = -24221.6665
$ws.Range("J65").Value = 3596.3333
$ws.Range("H65").Value = 3596.3333
$ws.Range("L65").Value = 17981.6665
$ws.Range("I122").Value = 3578.3076
$ws.Range("K122").Value = 10734.9228
$ws.Range("H122").Value = 3191.7896
$ws.Range("M122").Value = -8284.9228
$ws.Range("M132").Value = -2441
$ws.Range("H132").Value = 1959.909
$ws.Range("K132").Value = 4971
$ws.Range("N132").Value = -20027
$ws.Range("I132").Value = 1657
$ws.Range("L132").Value = 14967
$ws.Range("J132").Value = 4989
$ws.Range("M134").Value = -4355.000100000001
$ws.Range("I134").Value = 2296.6667
$ws.Range("K134").Value = 6890.000100000001
$ws.Range("H134").Value = 2278
$ws.Range("N140").Value = -50360
$ws.Range("J140").Value = 40000
$ws.Range("H140").Value = 40000
$ws.Range("L140").Value = 40000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 15686.15
$ws.Range("M4").Value = -46946.45
$ws.Range("H4").Value = 1443510.6
$ws.Range("K4").Value = 47058.45
$ws.Range("K92").Value = 751.5
$ws.Range("M92").Value = 496.5
$ws.Range("H92").Value = 270.2
$ws.Range("I92").Value = 250.5
$ws.Range("H97").Value = 1115.1428
$ws.Range("N97").Value = -4512.200000000001
$ws.Range("J97").Value = 1173.4
$ws.Range("L97").Value = 3520.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51372
$ws.Range("J62").Value = 50000
$ws.Range("N64").Value = -190767
$ws.Range("H64").Value = 190271
$ws.Range("J64").Value = 190271
$ws.Range("L64").Value = 190271
$ws.Range("N65").Value = -156864
$ws.Range("J65").Value = 50000
$ws.Range("H65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("J67").Value = 190271
$ws.Range("N67").Value = -191987
$ws.Range("L67").Value = 190271
$ws.Range("H67").Value = 190271
$ws.Range("I97").Value = 452.16666
$ws.Range("H97").Value = 480.0435
$ws.Range("K97").Value = 452.16666
$ws.Range("M97").Value = 43.83334000000002
$ws.Range("H126").Value = 1839
$ws.Range("I126").Value = 1839
$ws.Range("M126").Value = -3047
$ws.Range("K126").Value = 5517
$ws.Range("M132").Value = -3900.9095
$ws.Range("H132").Value = 2380.818
$ws.Range("K132").Value = 6430.9095
$ws.Range("N132").Value = -12914
$ws.Range("I132").Value = 2143.6365
$ws.Range("L132").Value = 7854
$ws.Range("J132").Value = 2618

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L22").Value = 1687.5
$ws.Range("K22").Value = 1469
$ws.Range("M22").Value = -1174
$ws.Range("I22").Value = 1469
$ws.Range("H22").Value = 1578.25
$ws.Range("J22").Value = 1687.5
$ws.Range("N22").Value = -2277.5
$ws.Range("J27").Value = 1687.5
$ws.Range("L27").Value = 1687.5
$ws.Range("K27").Value = 1469
$ws.Range("H27").Value = 1578.25
$ws.Range("I27").Value = 1469
$ws.Range("N27").Value = -1901.5
$ws.Range("M27").Value = -1362
$ws.Range("J55").Value = 1593.5
$ws.Range("L55").Value = 1593.5
$ws.Range("H55").Value = 1162.25
$ws.Range("N55").Value = -1939.5
$ws.Range("H61").Value = 1699.3334
$ws.Range("I61").Value = 1699.3334
$ws.Range("K61").Value = 1699.3334
$ws.Range("M61").Value = -1497.3334
$ws.Range("H82").Value = 1170.125
$ws.Range("N82").Value = -1233.8
$ws.Range("L82").Value = 511.8
$ws.Range("I82").Value = 2267.3333
$ws.Range("K82").Value = 2267.3333
$ws.Range("M82").Value = -1906.3333
$ws.Range("J82").Value = 511.8
$ws.Range("N85").Value = -3007.8
$ws.Range("M85").Value = -1019.3333
$ws.Range("L85").Value = 511.8
$ws.Range("J85").Value = 511.8
$ws.Range("H85").Value = 1170.125
$ws.Range("I85").Value = 2267.3333
$ws.Range("K85").Value = 2267.3333
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("N94").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N100").Value = -1582
$ws.Range("M100").Value = -125.6667
$ws.Range("I100").Value = 666.6667
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 666.6667
$ws.Range("L100").Value = 500
$ws.Range("H100").Value = 625
$ws.Range("H113").Value = 1699.3334
$ws.Range("M113").Value = 470.6666
$ws.Range("K113").Value = 1699.3334
$ws.Range("I113").Value = 1699.3334
$ws.Range("M132").Value = -8073.071599999999
$ws.Range("H132").Value = 3599.1177
$ws.Range("K132").Value = 10603.0716
$ws.Range("I132").Value = 3534.3572
$ws.Range("K136").Value = 5847
$ws.Range("I136").Value = 1949
$ws.Range("H136").Value = 1949
$ws.Range("M136").Value = -3297

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 99994.5
$ws.Range("J138").Value = 99994.5
$ws.Range("N138").Value = -110274.5
$ws.Range("L138").Value = 99994.5
